$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed "hot stock" cells for rows 2-21 (columns A/B/C)
$ws.Range("A2").Value = "首开股份"
$ws.Range("B2").Value = "上海建工"
$ws.Range("C2").Value = "上海建工"
$ws.Range("B3").Value = "卧龙电驱"
$ws.Range("C3").Value = "卧龙电驱"
$ws.Range("B4").Value = "山子高科"
$ws.Range("C4").Value = "吉视传媒"
$ws.Range("A5").Value = "上海建工"
$ws.Range("B5").Value = "首开股份"
$ws.Range("C5").Value = "山子高科"
$ws.Range("A6").Value = "立讯精密"
$ws.Range("C6").Value = "首开股份"
$ws.Range("A7").Value = "山子高科"
$ws.Range("B7").Value = "方正科技"
$ws.Range("C7").Value = "青山纸业"
$ws.Range("A8").Value = "青山纸业"
$ws.Range("B8").Value = "工业富联"
$ws.Range("C8").Value = "工业富联"
$ws.Range("A9").Value = "吉视传媒"
$ws.Range("C9").Value = "方正科技"
$ws.Range("A10").Value = "方正科技"
$ws.Range("B10").Value = "北方铜业"
$ws.Range("C10").Value = "华胜天成"
$ws.Range("A11").Value = "寒武纪-U"
$ws.Range("B11").Value = "立讯精密"
$ws.Range("C11").Value = "北方稀土"
$ws.Range("B12").Value = "北方稀土"
$ws.Range("C12").Value = "金发科技"
$ws.Range("A13").Value = "北方稀土"
$ws.Range("C13").Value = "利欧股份"
$ws.Range("A14").Value = "胜宏科技"
$ws.Range("C14").Value = "*ST东通"
$ws.Range("A15").Value = "金发科技"
$ws.Range("C15").Value = "先导智能"
$ws.Range("A16").Value = "赢合科技"
$ws.Range("C16").Value = "国轩高科"
$ws.Range("A17").Value = "国轩高科"
$ws.Range("C17").Value = "胜宏科技"
$ws.Range("A18").Value = "北方铜业"
$ws.Range("C18").Value = "立讯精密"
$ws.Range("A19").Value = "蜀道装备"
$ws.Range("C19").Value = "淳中科技"
$ws.Range("A20").Value = "厚普股份"
$ws.Range("C20").Value = "赢合科技"
$ws.Range("A21").Value = "天普股份"
$ws.Range("C21").Value = "三维通信"

# Column B is cleared out for rows 13-21 in the updated snapshot
$ws.Range("B13:B21").ClearContents()

